$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.337.20"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.578.23"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.491"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.48"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.82"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0896"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("D13").Value = "1.803.64"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "1.574.94"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.516"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "28.350.65"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0480"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.57%  "
$ws.Range("E32").Value = "  -3.27%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("D35").Value = "1.391.06"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.31%  "
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.59%  "
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.517"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.67%  "
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.785"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0456"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.73%  "
$ws.Range("E47").Value = "  -4.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "62.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = "1.715.12"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "41.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.20%  "
